# Update the attendance sheet so that each date row's "Invalid" (row 3 only)
# and "Absent" (rows 3-18) counts reflect the final values from the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 also gets its "Invalid" (column G) count updated.
$ws.Range("G3").Value = 1

# Column H ("Absent") goes from 0 to 1 for every data row, 3 through 18.
for ($row = 3; $row -le 18; $row++) {
    $ws.Cells.Item($row, 8).Value = 1
}
